$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 68.443746
$ws.Cells.Item(2, 8).Value = 205.331238
$ws.Cells.Item(2, 9).Value = 0.1596169534001499
$ws.Cells.Item(2, 10).Value = 0.1596169534001499
$ws.Cells.Item(2, 13).Value = 1.918906333333333
$ws.Cells.Item(2, 14).Value = 5.756718999999999
$ws.Cells.Item(2, 15).Value = 0.006524019162508824
$ws.Cells.Item(2, 16).Value = 0.006524019162508824
$ws.Cells.Item(2, 17).Value = 131.337137676458
$ws.Cells.Item(2, 18).Value = 1182.034239088122
$ws.Cells.Item(2, 19).Value = 0.001041344062643856
$ws.Cells.Item(2, 20).Value = 0.001041344062643856

# Row 3
$ws.Cells.Item(3, 7).Value = 68.443746
$ws.Cells.Item(3, 8).Value = 205.331238
$ws.Cells.Item(3, 9).Value = 0.1596169534001499
$ws.Cells.Item(3, 10).Value = 0.1596169534001499
$ws.Cells.Item(3, 13).Value = 181.2883913333334
$ws.Cells.Item(3, 14).Value = 543.865174
$ws.Cells.Item(3, 15).Value = 0.6163557430885885
$ws.Cells.Item(3, 16).Value = 0.6163557430885885
$ws.Cells.Item(3, 17).Value = 12408.05660916727
$ws.Cells.Item(3, 18).Value = 111672.5094825054
$ws.Cells.Item(3, 19).Value = 0.09838082592248601
$ws.Cells.Item(3, 20).Value = 0.09838082592248601

# Row 4
$ws.Cells.Item(4, 7).Value = 68.443746
$ws.Cells.Item(4, 8).Value = 205.331238
$ws.Cells.Item(4, 9).Value = 0.1596169534001499
$ws.Cells.Item(4, 10).Value = 0.1596169534001499
$ws.Cells.Item(4, 13).Value = 29.04767233333333
$ws.Cells.Item(4, 14).Value = 87.143017
$ws.Cells.Item(4, 15).Value = 0.09875811426384234
$ws.Cells.Item(4, 16).Value = 0.09875811426384236
$ws.Cells.Item(4, 17).Value = 1988.131507073894
$ws.Cells.Item(4, 18).Value = 17893.18356366505
$ws.Cells.Item(4, 19).Value = 0.0157634693223384
$ws.Cells.Item(4, 20).Value = 0.01576346932233841

# Row 5
$ws.Cells.Item(5, 7).Value = 68.443746
$ws.Cells.Item(5, 8).Value = 205.331238
$ws.Cells.Item(5, 9).Value = 0.1596169534001499
$ws.Cells.Item(5, 10).Value = 0.1596169534001499
$ws.Cells.Item(5, 13).Value = 81.87450533333333
$ws.Cells.Item(5, 14).Value = 245.623516
$ws.Cells.Item(5, 15).Value = 0.2783621234850603
$ws.Cells.Item(5, 16).Value = 0.2783621234850603
$ws.Cells.Item(5, 17).Value = 5603.797846910312
$ws.Cells.Item(5, 18).Value = 50434.18062219281
$ws.Cells.Item(5, 19).Value = 0.04443131409268165
$ws.Cells.Item(5, 20).Value = 0.04443131409268165

# Row 6
$ws.Cells.Item(6, 7).Value = 178.365814
$ws.Cells.Item(6, 8).Value = 535.097442
$ws.Cells.Item(6, 9).Value = 0.4159650732941736
$ws.Cells.Item(6, 10).Value = 0.4159650732941736
$ws.Cells.Item(6, 13).Value = 1.918906333333333
$ws.Cells.Item(6, 14).Value = 5.756718999999999
$ws.Cells.Item(6, 15).Value = 0.006524019162508824
$ws.Cells.Item(6, 16).Value = 0.006524019162508824
$ws.Cells.Item(6, 17).Value = 342.2672901347553
$ws.Cells.Item(6, 18).Value = 3080.405611212798
$ws.Cells.Item(6, 19).Value = 0.002713764109105576
$ws.Cells.Item(6, 20).Value = 0.002713764109105576

# Row 7
$ws.Cells.Item(7, 7).Value = 178.365814
$ws.Cells.Item(7, 8).Value = 535.097442
$ws.Cells.Item(7, 9).Value = 0.4159650732941736
$ws.Cells.Item(7, 10).Value = 0.4159650732941736
$ws.Cells.Item(7, 13).Value = 181.2883913333334
$ws.Cells.Item(7, 14).Value = 543.865174
$ws.Cells.Item(7, 15).Value = 0.6163557430885885
$ws.Cells.Item(7, 16).Value = 0.6163557430885885
$ws.Cells.Item(7, 17).Value = 32335.65148892055
$ws.Cells.Item(7, 18).Value = 291020.8634002849
$ws.Cells.Item(7, 19).Value = 0.2563824618491296
$ws.Cells.Item(7, 20).Value = 0.2563824618491296

# Row 8
$ws.Cells.Item(8, 7).Value = 178.365814
$ws.Cells.Item(8, 8).Value = 535.097442
$ws.Cells.Item(8, 9).Value = 0.4159650732941736
$ws.Cells.Item(8, 10).Value = 0.4159650732941736
$ws.Cells.Item(8, 13).Value = 29.04767233333333
$ws.Cells.Item(8, 14).Value = 87.143017
$ws.Cells.Item(8, 15).Value = 0.09875811426384234
$ws.Cells.Item(8, 16).Value = 0.09875811426384236
$ws.Cells.Item(8, 17).Value = 5181.11172054028
$ws.Cells.Item(8, 18).Value = 46630.00548486252
$ws.Cells.Item(8, 19).Value = 0.04107992623815356
$ws.Cells.Item(8, 20).Value = 0.04107992623815356

# Row 9
$ws.Cells.Item(9, 7).Value = 178.365814
$ws.Cells.Item(9, 8).Value = 535.097442
$ws.Cells.Item(9, 9).Value = 0.4159650732941736
$ws.Cells.Item(9, 10).Value = 0.4159650732941736
$ws.Cells.Item(9, 13).Value = 81.87450533333333
$ws.Cells.Item(9, 14).Value = 245.623516
$ws.Cells.Item(9, 15).Value = 0.2783621234850603
$ws.Cells.Item(9, 16).Value = 0.2783621234850603
$ws.Cells.Item(9, 17).Value = 14603.61278962734
$ws.Cells.Item(9, 18).Value = 131432.5151066461
$ws.Cells.Item(9, 19).Value = 0.1157889210977849
$ws.Cells.Item(9, 20).Value = 0.1157889210977849

# Row 10
$ws.Cells.Item(10, 7).Value = 88.88346833333333
$ws.Cells.Item(10, 8).Value = 266.650405
$ws.Cells.Item(10, 9).Value = 0.2072842188241036
$ws.Cells.Item(10, 10).Value = 0.2072842188241036
$ws.Cells.Item(10, 13).Value = 1.918906333333333
$ws.Cells.Item(10, 14).Value = 5.756718999999999
$ws.Cells.Item(10, 15).Value = 0.006524019162508824
$ws.Cells.Item(10, 16).Value = 0.006524019162508824
$ws.Cells.Item(10, 17).Value = 170.5590503134661
$ws.Cells.Item(10, 18).Value = 1535.031452821195
$ws.Cells.Item(10, 19).Value = 0.001352326215694124
$ws.Cells.Item(10, 20).Value = 0.001352326215694124

# Row 11
$ws.Cells.Item(11, 7).Value = 88.88346833333333
$ws.Cells.Item(11, 8).Value = 266.650405
$ws.Cells.Item(11, 9).Value = 0.2072842188241036
$ws.Cells.Item(11, 10).Value = 0.2072842188241036
$ws.Cells.Item(11, 13).Value = 181.2883913333334
$ws.Cells.Item(11, 14).Value = 543.865174
$ws.Cells.Item(11, 15).Value = 0.6163557430885885
$ws.Cells.Item(11, 16).Value = 0.6163557430885885
$ws.Cells.Item(11, 17).Value = 16113.54099027728
$ws.Cells.Item(11, 18).Value = 145021.8689124955
$ws.Cells.Item(11, 19).Value = 0.127760818723868
$ws.Cells.Item(11, 20).Value = 0.127760818723868

# Row 12
$ws.Cells.Item(12, 7).Value = 88.88346833333333
$ws.Cells.Item(12, 8).Value = 266.650405
$ws.Cells.Item(12, 9).Value = 0.2072842188241036
$ws.Cells.Item(12, 10).Value = 0.2072842188241036
$ws.Cells.Item(12, 13).Value = 29.04767233333333
$ws.Cells.Item(12, 14).Value = 87.143017
$ws.Cells.Item(12, 15).Value = 0.09875811426384234
$ws.Cells.Item(12, 16).Value = 0.09875811426384236
$ws.Cells.Item(12, 17).Value = 2581.857863996876
$ws.Cells.Item(12, 18).Value = 23236.72077597188
$ws.Cells.Item(12, 19).Value = 0.02047099856772212
$ws.Cells.Item(12, 20).Value = 0.02047099856772212

# Row 13
$ws.Cells.Item(13, 7).Value = 88.88346833333333
$ws.Cells.Item(13, 8).Value = 266.650405
$ws.Cells.Item(13, 9).Value = 0.2072842188241036
$ws.Cells.Item(13, 10).Value = 0.2072842188241036
$ws.Cells.Item(13, 13).Value = 81.87450533333333
$ws.Cells.Item(13, 14).Value = 245.623516
$ws.Cells.Item(13, 15).Value = 0.2783621234850603
$ws.Cells.Item(13, 16).Value = 0.2783621234850603
$ws.Cells.Item(13, 17).Value = 7277.290002102663
$ws.Cells.Item(13, 18).Value = 65495.61001892397
$ws.Cells.Item(13, 19).Value = 0.05770007531681938
$ws.Cells.Item(13, 20).Value = 0.05770007531681938

# Row 14
$ws.Cells.Item(14, 7).Value = 93.106949
$ws.Cells.Item(14, 8).Value = 279.320847
$ws.Cells.Item(14, 9).Value = 0.2171337544815728
$ws.Cells.Item(14, 10).Value = 0.2171337544815728
$ws.Cells.Item(14, 13).Value = 1.918906333333333
$ws.Cells.Item(14, 14).Value = 5.756718999999999
$ws.Cells.Item(14, 15).Value = 0.006524019162508824
$ws.Cells.Item(14, 16).Value = 0.006524019162508824
$ws.Cells.Item(14, 17).Value = 178.6635141134437
$ws.Cells.Item(14, 18).Value = 1607.971627020993
$ws.Cells.Item(14, 19).Value = 0.001416584775065267
$ws.Cells.Item(14, 20).Value = 0.001416584775065267

# Row 15
$ws.Cells.Item(15, 7).Value = 93.106949
$ws.Cells.Item(15, 8).Value = 279.320847
$ws.Cells.Item(15, 9).Value = 0.2171337544815728
$ws.Cells.Item(15, 10).Value = 0.2171337544815728
$ws.Cells.Item(15, 13).Value = 181.2883913333334
$ws.Cells.Item(15, 14).Value = 543.865174
$ws.Cells.Item(15, 15).Value = 0.6163557430885885
$ws.Cells.Item(15, 16).Value = 0.6163557430885885
$ws.Cells.Item(15, 17).Value = 16879.20900616471
$ws.Cells.Item(15, 18).Value = 151912.8810554824
$ws.Cells.Item(15, 19).Value = 0.133831636593105
$ws.Cells.Item(15, 20).Value = 0.133831636593105

# Row 16
$ws.Cells.Item(16, 7).Value = 93.106949
$ws.Cells.Item(16, 8).Value = 279.320847
$ws.Cells.Item(16, 9).Value = 0.2171337544815728
$ws.Cells.Item(16, 10).Value = 0.2171337544815728
$ws.Cells.Item(16, 13).Value = 29.04767233333333
$ws.Cells.Item(16, 14).Value = 87.143017
$ws.Cells.Item(16, 15).Value = 0.09875811426384234
$ws.Cells.Item(16, 16).Value = 0.09875811426384236
$ws.Cells.Item(16, 17).Value = 2704.540146508378
$ws.Cells.Item(16, 18).Value = 24340.8613185754
$ws.Cells.Item(16, 19).Value = 0.02144372013562826
$ws.Cells.Item(16, 20).Value = 0.02144372013562827

# Row 17
$ws.Cells.Item(17, 7).Value = 93.106949
$ws.Cells.Item(17, 8).Value = 279.320847
$ws.Cells.Item(17, 9).Value = 0.2171337544815728
$ws.Cells.Item(17, 10).Value = 0.2171337544815728
$ws.Cells.Item(17, 13).Value = 81.87450533333333
$ws.Cells.Item(17, 14).Value = 245.623516
$ws.Cells.Item(17, 15).Value = 0.2783621234850603
$ws.Cells.Item(17, 16).Value = 0.2783621234850603
$ws.Cells.Item(17, 17).Value = 7623.085392470894
$ws.Cells.Item(17, 18).Value = 68607.76853223806
$ws.Cells.Item(17, 19).Value = 0.06044181297777435
$ws.Cells.Item(17, 20).Value = 0.06044181297777435
